$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto-price snapshot refresh: update Price (D) and Volume(1h) (E)
# columns for the rows whose source data changed in this run.
# D-column values are plain text (e.g. "28.639.86", "1.004") — for the ones that
# look like a plain decimal number (single "." or none) we force the cell to Text
# format first so Excel does not silently reinterpret the string as a Number.

$ws.Range("D2").Value = "28.639.86"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "1.868.73"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.26"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3881"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07871"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9750"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.848.95"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.003"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.701"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06963"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.25"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001002"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "28.641.24"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.272"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.106"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "2.095.48"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.67"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.23"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.857"
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.986"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.33"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09321"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9179"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.278"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.319"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05793"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.151"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02076"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.686"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5627"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1784"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.765"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07223"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.76"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.151"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.126"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.78"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.409"
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  +0.17%  "
